$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new comment to D20 and D23 for students who didn't attend the lab
# and didn't show the TA the result.
$comment = "Didn't attend the lab, didn't show TA the result. "
$ws.Range("D20").Value = $comment
$ws.Range("D23").Value = $comment

# Update the saved view state of the sheet (scroll position / selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
